$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 151; this shifts the existing rows 151..169 down to 152..170
# (keeping all of their original values/styles intact).
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record's data.
$ws.Cells.Item(151,1).Value = 10
$ws.Cells.Item(151,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(151,3).Value = "La Araucanía"
$ws.Cells.Item(151,4).Value = 44449
$ws.Cells.Item(151,5).Value = 9
$ws.Cells.Item(151,6).Value = 100112017
$ws.Cells.Item(151,7).Value = "Apio"
$ws.Cells.Item(151,8).Value = "Americana (o)"
$ws.Cells.Item(151,9).Value = "Primera"
$ws.Cells.Item(151,10).Value = 95
$ws.Cells.Item(151,11).Value = 10000
$ws.Cells.Item(151,12).Value = 10000
$ws.Cells.Item(151,13).Value = 10000
$ws.Cells.Item(151,14).Value = "`$/docena de matas"
$ws.Cells.Item(151,15).Value = "Provincia del Elquí"
$ws.Cells.Item(151,16).Value = 1667
$ws.Cells.Item(151,17).Value = 6
$ws.Cells.Item(151,18).Value = "Hortaliza"
